$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.609.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.85%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.867.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.27%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'235.15"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.08%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4704"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.37%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.2767"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.48%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.06377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.13%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'17.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +10.48%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'1.864.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.03%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +0.28%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'4.983"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.12%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'85.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.79%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.6360"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.66%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'30.581.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.94%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'241.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.52%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.9998"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.01%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'12.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.28%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.000007379"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.08%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.22%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'4.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.88%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'6.037"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.48%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'9.395"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.12%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'165.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.57%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'18.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.78%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  +1.74%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'0.1023"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.18%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'1.379"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.16%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'4.103"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.62%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'3.868"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.46%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.04938"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.51%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.153"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.53%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.7094"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.99%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'2.705"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.25%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.01910"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.53%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'2.689"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.11%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.8816"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.38%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'2.000"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.91%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'105.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.03%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'1.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.01%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.4113"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.18%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'5.556"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.16%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'7.417"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +5.05%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'62.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.62%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.1231"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.96%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'33.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.11%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'8.647"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.80%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05574"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.44%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.381"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.62%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.3719"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.57%  "
$ws.Range("E51").Style = "Normal"
